$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 49.88947433333334
$ws.Range("H2").Value = 149.668423
$ws.Range("I2").Value = 0.2324880572195875
$ws.Range("J2").Value = 0.2324880572195874
$ws.Range("M2").Value = 7.938978333333334
$ws.Range("N2").Value = 23.816935
$ws.Range("O2").Value = 0.2711338618634719
$ws.Range("P2").Value = 0.2711338618634719
$ws.Range("Q2").Value = 396.0714557937229
$ws.Range("R2").Value = 3564.643102143506
$ws.Range("S2").Value = 0.06303538479108257
$ws.Range("T2").Value = 0.06303538479108256

$ws.Range("G3").Value = 49.88947433333334
$ws.Range("H3").Value = 149.668423
$ws.Range("I3").Value = 0.2324880572195875
$ws.Range("J3").Value = 0.2324880572195874
$ws.Range("M3").Value = 7.621805666666667
$ws.Range("O3").Value = 0.2603017060897501
$ws.Range("P3").Value = 0.2603017060897501
$ws.Range("Q3").Value = 380.2478781808213
$ws.Range("R3").Value = 3422.230903627391
$ws.Range("S3").Value = 0.06051703793975005
$ws.Range("T3").Value = 0.06051703793975004

$ws.Range("G4").Value = 49.88947433333334
$ws.Range("H4").Value = 149.668423
$ws.Range("I4").Value = 0.2324880572195875
$ws.Range("J4").Value = 0.2324880572195874
$ws.Range("M4").Value = 2.437389666666667
$ws.Range("N4").Value = 7.312169
$ws.Range("O4").Value = 0.08324230718891248
$ws.Range("P4").Value = 0.08324230718891248
$ws.Range("Q4").Value = 121.6000892154986
$ws.Range("R4").Value = 1094.400802939487
$ws.Range("S4").Value = 0.01935284227682636
$ws.Range("T4").Value = 0.01935284227682636

$ws.Range("G5").Value = 49.88947433333334
$ws.Range("H5").Value = 149.668423
$ws.Range("I5").Value = 0.2324880572195875
$ws.Range("J5").Value = 0.2324880572195874
$ws.Range("M5").Value = 11.282486
$ws.Range("N5").Value = 33.847458
$ws.Range("O5").Value = 0.3853221248578654
$ws.Range("P5").Value = 0.3853221248578655
$ws.Range("Q5").Value = 562.8772957131928
$ws.Range("R5").Value = 5065.895661418735
$ws.Range("S5").Value = 0.08958279221192844
$ws.Range("T5").Value = 0.08958279221192844

$ws.Range("I6").Value = 0.295249080025651
$ws.Range("J6").Value = 0.295249080025651
$ws.Range("M6").Value = 7.938978333333334
$ws.Range("N6").Value = 23.816935
$ws.Range("O6").Value = 0.2711338618634719
$ws.Range("P6").Value = 0.2711338618634719
$ws.Range("Q6").Value = 502.9924304329584
$ws.Range("R6").Value = 4526.931873896626
$ws.Range("S6").Value = 0.08005202327899202
$ws.Range("T6").Value = 0.08005202327899202

$ws.Range("I7").Value = 0.295249080025651
$ws.Range("J7").Value = 0.295249080025651
$ws.Range("M7").Value = 7.621805666666667
$ws.Range("O7").Value = 0.2603017060897501
$ws.Range("P7").Value = 0.2603017060897501
$ws.Range("Q7").Value = 482.8972187098417
$ws.Range("R7").Value = 4346.074968388575
$ws.Range("S7").Value = 0.07685383925210611
$ws.Range("T7").Value = 0.07685383925210611

$ws.Range("I8").Value = 0.295249080025651
$ws.Range("J8").Value = 0.295249080025651
$ws.Range("M8").Value = 2.437389666666667
$ws.Range("N8").Value = 7.312169
$ws.Range("O8").Value = 0.08324230718891248
$ws.Range("P8").Value = 0.08324230718891248
$ws.Range("Q8").Value = 154.4264892626417
$ws.Range("R8").Value = 1389.838403363775
$ws.Range("S8").Value = 0.02457721461673905
$ws.Range("T8").Value = 0.02457721461673905

$ws.Range("I9").Value = 0.295249080025651
$ws.Range("J9").Value = 0.295249080025651
$ws.Range("M9").Value = 11.282486
$ws.Range("N9").Value = 33.847458
$ws.Range("O9").Value = 0.3853221248578654
$ws.Range("P9").Value = 0.3853221248578655
$ws.Range("Q9").Value = 714.8281323099501
$ws.Range("R9").Value = 6433.453190789551
$ws.Range("S9").Value = 0.1137660028778138
$ws.Range("T9").Value = 0.1137660028778138

$ws.Range("G10").Value = 52.37451933333333
$ws.Range("H10").Value = 157.123558
$ws.Range("I10").Value = 0.2440685216737345
$ws.Range("J10").Value = 0.2440685216737345
$ws.Range("M10").Value = 7.938978333333334
$ws.Range("N10").Value = 23.816935
$ws.Range("O10").Value = 0.2711338618634719
$ws.Range("P10").Value = 0.2711338618634719
$ws.Range("Q10").Value = 415.8001742060811
$ws.Range("R10").Value = 3742.20156785473
$ws.Range("S10").Value = 0.06617524084070812
$ws.Range("T10").Value = 0.06617524084070812

$ws.Range("G11").Value = 52.37451933333333
$ws.Range("H11").Value = 157.123558
$ws.Range("I11").Value = 0.2440685216737345
$ws.Range("J11").Value = 0.2440685216737345
$ws.Range("M11").Value = 7.621805666666667
$ws.Range("O11").Value = 0.2603017060897501
$ws.Range("P11").Value = 0.2603017060897501
$ws.Range("Q11").Value = 399.1884082437429
$ws.Range("R11").Value = 3592.695674193686
$ws.Range("S11").Value = 0.06353145259447623
$ws.Range("T11").Value = 0.06353145259447623

$ws.Range("G12").Value = 52.37451933333333
$ws.Range("H12").Value = 157.123558
$ws.Range("I12").Value = 0.2440685216737345
$ws.Range("J12").Value = 0.2440685216737345
$ws.Range("M12").Value = 2.437389666666667
$ws.Range("N12").Value = 7.312169
$ws.Range("O12").Value = 0.08324230718891248
$ws.Range("P12").Value = 0.08324230718891248
$ws.Range("Q12").Value = 127.6571122197002
$ws.Range("R12").Value = 1148.914009977302
$ws.Range("S12").Value = 0.02031682685630875
$ws.Range("T12").Value = 0.02031682685630875

$ws.Range("G13").Value = 52.37451933333333
$ws.Range("H13").Value = 157.123558
$ws.Range("I13").Value = 0.2440685216737345
$ws.Range("J13").Value = 0.2440685216737345
$ws.Range("M13").Value = 11.282486
$ws.Range("N13").Value = 33.847458
$ws.Range("O13").Value = 0.3853221248578654
$ws.Range("P13").Value = 0.3853221248578655
$ws.Range("Q13").Value = 590.9147811350626
$ws.Range("R13").Value = 5318.233030215564
$ws.Range("S13").Value = 0.09404500138224135
$ws.Range("T13").Value = 0.09404500138224137

$ws.Range("G14").Value = 48.96808833333333
$ws.Range("H14").Value = 146.904265
$ws.Range("I14").Value = 0.2281943410810271
$ws.Range("J14").Value = 0.228194341081027
$ws.Range("M14").Value = 7.938978333333334
$ws.Range("N14").Value = 23.816935
$ws.Range("O14").Value = 0.2711338618634719
$ws.Range("P14").Value = 0.2711338618634719
$ws.Range("Q14").Value = 388.7565923030862
$ws.Range("R14").Value = 3498.809330727775
$ws.Range("S14").Value = 0.06187121295268918
$ws.Range("T14").Value = 0.06187121295268917

$ws.Range("G15").Value = 48.96808833333333
$ws.Range("H15").Value = 146.904265
$ws.Range("I15").Value = 0.2281943410810271
$ws.Range("J15").Value = 0.228194341081027
$ws.Range("M15").Value = 7.621805666666667
$ws.Range("O15").Value = 0.2603017060897501
$ws.Range("P15").Value = 0.2603017060897501
$ws.Range("Q15").Value = 373.2252531448339
$ws.Range("R15").Value = 3359.027278303505
$ws.Range("S15").Value = 0.05939937630341768
$ws.Range("T15").Value = 0.05939937630341768

$ws.Range("G16").Value = 48.96808833333333
$ws.Range("H16").Value = 146.904265
$ws.Range("I16").Value = 0.2281943410810271
$ws.Range("J16").Value = 0.228194341081027
$ws.Range("M16").Value = 2.437389666666667
$ws.Range("N16").Value = 7.312169
$ws.Range("O16").Value = 0.08324230718891248
$ws.Range("P16").Value = 0.08324230718891248
$ws.Range("Q16").Value = 119.3543125000872
$ws.Range("R16").Value = 1074.188812500785
$ws.Range("S16").Value = 0.01899542343903833
$ws.Range("T16").Value = 0.01899542343903832

$ws.Range("G17").Value = 48.96808833333333
$ws.Range("H17").Value = 146.904265
$ws.Range("I17").Value = 0.2281943410810271
$ws.Range("J17").Value = 0.2281943410810271
$ws.Range("M17").Value = 11.282486
$ws.Range("N17").Value = 33.847458
$ws.Range("O17").Value = 0.3853221248578654
$ws.Range("P17").Value = 0.3853221248578655
$ws.Range("Q17").Value = 552.4817710675967
$ws.Range("R17").Value = 4972.335939608371
$ws.Range("S17").Value = 0.08792832838588184
$ws.Range("T17").Value = 0.08792832838588184
